# Update the single data cell on the first sheet: the hyperlink-style text
# that previously pointed at forex4you is replaced with the litefinance URL
# (commit: "fixed lifefinance crash with name").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "https://my.litefinance.org/ru/traders/trades-history?id=2187036"

# The author's last selection on the sheet moved from A8 to E13 before saving.
$null = $ws.Range("E13").Select()
